$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Unprotect()
$ws.Range('A1').Value = '!!!ObjTables objTablesVersion=''0.0.8'' date=''2020-03-09 15:31:39'''
$ws.Range('A2').Value = '!!ObjTables type=''Data'' id=''Compartment'' name=''Compartment'' date=''2020-03-09 15:31:39'' objTablesVersion=''0.0.8'' tableFormat=''row'''
$ws.Protect()

$ws = $wb.Worksheets.Item(2)
$ws.Unprotect()
$ws.Range('A1').Value = '!!ObjTables type=''Data'' id=''Compound'' name=''Compound'' date=''2020-03-09 15:31:39'' objTablesVersion=''0.0.8'' tableFormat=''row'''
$ws.Protect()

$ws = $wb.Worksheets.Item(3)
$ws.Unprotect()
$ws.Range('A1').Value = '!!ObjTables type=''Data'' id=''Definition'' name=''Definition'' date=''2020-03-09 15:31:39'' objTablesVersion=''0.0.8'' tableFormat=''row'''
$ws.Protect()

$ws = $wb.Worksheets.Item(4)
$ws.Unprotect()
$ws.Range('A1').Value = '!!ObjTables type=''Data'' id=''Enzyme'' name=''Enzyme'' date=''2020-03-09 15:31:39'' objTablesVersion=''0.0.8'' tableFormat=''row'''
$ws.Protect()

$ws = $wb.Worksheets.Item(5)
$ws.Unprotect()
$ws.Range('A1').Value = '!!ObjTables type=''Data'' id=''FbcObjective'' name=''FbcObjective'' date=''2020-03-09 15:31:39'' objTablesVersion=''0.0.8'' tableFormat=''row'''
$ws.Protect()

$ws = $wb.Worksheets.Item(6)
$ws.Unprotect()
$ws.Range('A1').Value = '!!ObjTables type=''Data'' id=''Gene'' name=''Gene'' date=''2020-03-09 15:31:39'' objTablesVersion=''0.0.8'' tableFormat=''row'''
$ws.Protect()

$ws = $wb.Worksheets.Item(7)
$ws.Unprotect()
$ws.Range('A1').Value = '!!ObjTables type=''Data'' id=''Layout'' name=''Layout'' date=''2020-03-09 15:31:39'' objTablesVersion=''0.0.8'' tableFormat=''row'''
$ws.Protect()

$ws = $wb.Worksheets.Item(8)
$ws.Unprotect()
$ws.Range('A1').Value = '!!ObjTables type=''Data'' id=''Measurement'' name=''Measurement'' date=''2020-03-09 15:31:39'' objTablesVersion=''0.0.8'' tableFormat=''row'''
$ws.Protect()

$ws = $wb.Worksheets.Item(9)
$ws.Unprotect()
$ws.Range('A1').Value = '!!ObjTables type=''Data'' id=''PbConfig'' name=''PbConfig'' date=''2020-03-09 15:31:39'' objTablesVersion=''0.0.8'' tableFormat=''row'''
$ws.Protect()

$ws = $wb.Worksheets.Item(10)
$ws.Unprotect()
$ws.Range('A1').Value = '!!ObjTables type=''Data'' id=''Position'' name=''Position'' date=''2020-03-09 15:31:39'' objTablesVersion=''0.0.8'' tableFormat=''row'''
$ws.Protect()

$ws = $wb.Worksheets.Item(11)
$ws.Unprotect()
$ws.Range('A1').Value = '!!ObjTables type=''Data'' id=''Protein'' name=''Protein'' date=''2020-03-09 15:31:39'' objTablesVersion=''0.0.8'' tableFormat=''row'''
$ws.Protect()

$ws = $wb.Worksheets.Item(12)
$ws.Unprotect()
$ws.Range('A1').Value = '!!ObjTables type=''Data'' id=''Quantity'' name=''Quantity'' date=''2020-03-09 15:31:39'' objTablesVersion=''0.0.8'' tableFormat=''row'''
$ws.Protect()

$ws = $wb.Worksheets.Item(13)
$ws.Unprotect()
$ws.Range('A1').Value = '!!ObjTables type=''Data'' id=''QuantityInfo'' name=''QuantityInfo'' date=''2020-03-09 15:31:39'' objTablesVersion=''0.0.8'' tableFormat=''row'''
$ws.Protect()

$ws = $wb.Worksheets.Item(14)
$ws.Unprotect()
$ws.Range('A1').Value = '!!ObjTables type=''Data'' id=''QuantityMatrix'' name=''QuantityMatrix'' date=''2020-03-09 15:31:39'' objTablesVersion=''0.0.8'' tableFormat=''row'''
$ws.Protect()

$ws = $wb.Worksheets.Item(15)
$ws.Unprotect()
$ws.Range('A1').Value = '!!ObjTables type=''Data'' id=''Reaction'' name=''Reaction'' date=''2020-03-09 15:31:39'' objTablesVersion=''0.0.8'' tableFormat=''row'''
$ws.Protect()

$ws = $wb.Worksheets.Item(16)
$ws.Unprotect()
$ws.Range('A1').Value = '!!ObjTables type=''Data'' id=''ReactionStoichiometry'' name=''ReactionStoichiometry'' date=''2020-03-09 15:31:40'' objTablesVersion=''0.0.8'' tableFormat=''row'''
$ws.Protect()

$ws = $wb.Worksheets.Item(17)
$ws.Unprotect()
$ws.Range('A1').Value = '!!ObjTables type=''Data'' id=''Regulator'' name=''Regulator'' date=''2020-03-09 15:31:40'' objTablesVersion=''0.0.8'' tableFormat=''row'''
$ws.Protect()

$ws = $wb.Worksheets.Item(18)
$ws.Unprotect()
$ws.Range('A1').Value = '!!ObjTables type=''Data'' id=''Relation'' name=''Relation'' date=''2020-03-09 15:31:40'' objTablesVersion=''0.0.8'' tableFormat=''row'''
$ws.Protect()

$ws = $wb.Worksheets.Item(19)
$ws.Unprotect()
$ws.Range('A1').Value = '!!ObjTables type=''Data'' id=''Relationship'' name=''Relationship'' date=''2020-03-09 15:31:40'' objTablesVersion=''0.0.8'' tableFormat=''row'''
$ws.Protect()

$ws = $wb.Worksheets.Item(20)
$ws.Unprotect()
$ws.Range('A1').Value = '!!ObjTables type=''Data'' id=''SparseMatrix'' name=''SparseMatrix'' date=''2020-03-09 15:31:40'' objTablesVersion=''0.0.8'' tableFormat=''row'''
$ws.Protect()

$ws = $wb.Worksheets.Item(21)
$ws.Unprotect()
$ws.Range('A1').Value = '!!ObjTables type=''Data'' id=''SparseMatrixColumn'' name=''SparseMatrixColumn'' date=''2020-03-09 15:31:40'' objTablesVersion=''0.0.8'' tableFormat=''row'''
$ws.Protect()

$ws = $wb.Worksheets.Item(22)
$ws.Unprotect()
$ws.Range('A1').Value = '!!ObjTables type=''Data'' id=''SparseMatrixOrdered'' name=''SparseMatrixOrdered'' date=''2020-03-09 15:31:40'' objTablesVersion=''0.0.8'' tableFormat=''row'''
$ws.Protect()

$ws = $wb.Worksheets.Item(23)
$ws.Unprotect()
$ws.Range('A1').Value = '!!ObjTables type=''Data'' id=''SparseMatrixRow'' name=''SparseMatrixRow'' date=''2020-03-09 15:31:40'' objTablesVersion=''0.0.8'' tableFormat=''row'''
$ws.Protect()

$ws = $wb.Worksheets.Item(24)
$ws.Unprotect()
$ws.Range('A1').Value = '!!ObjTables type=''Data'' id=''StoichiometricMatrix'' name=''StoichiometricMatrix'' date=''2020-03-09 15:31:40'' objTablesVersion=''0.0.8'' tableFormat=''row'''
$ws.Protect()

$ws = $wb.Worksheets.Item(25)
$ws.Unprotect()
$ws.Range('A1').Value = '!!ObjTables type=''Data'' id=''rxnconContingencyList'' name=''rxnconContingencyList'' date=''2020-03-09 15:31:40'' objTablesVersion=''0.0.8'' tableFormat=''row'''
$ws.Protect()

$ws = $wb.Worksheets.Item(26)
$ws.Unprotect()
$ws.Range('A1').Value = '!!ObjTables type=''Data'' id=''rxnconReactionList'' name=''rxnconReactionList'' date=''2020-03-09 15:31:40'' objTablesVersion=''0.0.8'' tableFormat=''row'''
$ws.Protect()
